# Updates the cryptos list (prices and 1h volume deltas) plus two
# row re-orderings (rows 14<->15 and rows 35<->36) to reflect the
# latest scrape, per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 14 and 15 swapped places (Avalanche now ranks above WrappedliquidstakedEther2.0) ---
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'27.58"
$ws.Range("E14").Value = "  -2.56%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.151.66"
$ws.Range("E15").Value = "  -1.33%  "

# --- Rows 35 and 36 swapped places (Kaspa now ranks above FirstDigitalUSD) ---
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.128"
$ws.Range("E35").Value = "  -3.08%  "

$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.12%  "

# --- Price (D) / Volume 1h change (E) refresh for the remaining rows ---
$ws.Range("D2").Value = "67.133.38"
$ws.Range("E2").Value = "  -2.06%  "
$ws.Range("D3").Value = "2.665.35"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'595.26"
$ws.Range("E5").Value = "  -0.57%  "
$ws.Range("D6").Value = "'164.32"
$ws.Range("E6").Value = "  +2.76%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.543"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").Value = "2.663.44"
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("D10").Value = "'0.140"
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").Value = "'0.355"
$ws.Range("E12").Value = "  -1.34%  "
$ws.Range("D13").Value = "'5.18"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("D16").Value = "'0.0000182"
$ws.Range("E16").Value = "  -3.18%  "
$ws.Range("D17").Value = "67.035.34"
$ws.Range("E17").Value = "  -2.14%  "
$ws.Range("D18").Value = "2.633.16"
$ws.Range("E18").Value = "  -2.36%  "
$ws.Range("D19").Value = "'11.63"
$ws.Range("E19").Value = "  -2.81%  "
$ws.Range("D20").Value = "'360.43"
$ws.Range("E20").Value = "  -1.67%  "
$ws.Range("D21").Value = "'7.49"
$ws.Range("E21").Value = "  -2.43%  "
$ws.Range("D22").Value = "'4.36"
$ws.Range("E22").Value = "  -4.27%  "
$ws.Range("D23").Value = "'4.78"
$ws.Range("E23").Value = "  -2.30%  "
$ws.Range("E24").Value = "  -5.15%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'71.05"
$ws.Range("E26").Value = "  -4.58%  "
$ws.Range("D27").Value = "'10.06"
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("E30").Value = "  -3.11%  "
$ws.Range("D31").Value = "'549.76"
$ws.Range("E31").Value = "  -4.02%  "
$ws.Range("D32").Value = "'7.94"
$ws.Range("E32").Value = "  -3.55%  "
$ws.Range("D33").Value = "'1.38"
$ws.Range("E33").Value = "  -5.62%  "
$ws.Range("E34").Value = "  -1.71%  "
$ws.Range("D37").Value = "'1.56"
$ws.Range("E37").Value = "  -5.33%  "
$ws.Range("D38").Value = "'19.44"
$ws.Range("E38").Value = "  -2.69%  "
$ws.Range("D39").Value = "'154.20"
$ws.Range("E39").Value = "  -4.23%  "
$ws.Range("D40").Value = "'0.371"
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("D41").Value = "'5.26"
$ws.Range("E41").Value = "  -3.02%  "
$ws.Range("E42").Value = "  -4.81%  "
$ws.Range("D43").Value = "'17.90"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "'2.50"
$ws.Range("E45").Value = "  -5.47%  "
$ws.Range("D46").Value = "'40.18"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").Value = "0.0₆0295"
$ws.Range("E47").Value = "  -6.89%  "
$ws.Range("D48").Value = "'0.584"
$ws.Range("E48").Value = "  -2.93%  "
$ws.Range("D49").Value = "'152.33"
$ws.Range("E49").Value = "  -3.85%  "
$ws.Range("D50").Value = "'3.81"
$ws.Range("E50").Value = "  -3.74%  "
$ws.Range("E51").Value = "  -3.84%  "
